# Turn the footer "Copie : ..." line from a hard-coded list of three
# fixed recipients (EDF OA/SEI/EDM + DREAL + CRE) into a templated loop
# over an "enCopie(s)" array, joining each entry with " ; ".
#
# Copie : {#isEDFOA}EDF OA{/isEDFOA}{#isEDFSEI}EDF SEI{/isEDFSEI}{#isEDM}EDM{/isEDM} ; {#dreal}DREAL {dreal}{/dreal}{^dreal}DREAL concernée{/dreal} ; CRE
#   ->
# Copie : {#enCopies} {.} ; {/enCopies}

$d = $word.ActiveDocument

$nbsp = [char]0x00A0

$old = "Copie" + $nbsp + ": {#isEDFOA}EDF OA{/isEDFOA}{#isEDFSEI}EDF SEI{/isEDFSEI}{#isEDM}EDM{/isEDM}" + `
       $nbsp + "; {#dreal}DREAL {dreal}{/dreal}{^dreal}DREAL concernée{/dreal} ; CRE"

$new = "Copie" + $nbsp + ": {#enCopies} {.} ; {/enCopies}"

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                                  $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find the 'Copie : ...' footer text to replace"
}
